$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header values: A1 keeps "EmpleadoNo", B1 changes from "Empresa Id" to "Entidad Id"
$ws.Range("A1").Value = "EmpleadoNo"
$ws.Range("B1").Value = "Entidad Id"

# Move the active selection to D3
$ws.Range("D3").Select()
